# Update "want to go" counts (column F) across the 展览 / 演出 / 全部类型 sheets,
# matching the new data snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map: worksheet name -> list of (row, newValue) pairs for column F.
# NOTE: a leading unary comma (,@( ... )) is required for single-pair lists
# so the outer wrapper array is not unrolled away by PowerShell.
$updates = @{
    "展览" = @(
        @(2, 1688),
        @(3, 9160),
        @(4, 116),
        @(5, 509),
        @(6, 712),
        @(7, 1380),
        @(8, 207),
        @(10, 101),
        @(11, 5946),
        @(13, 390),
        @(15, 4589),
        @(19, 34),
        @(20, 341),
        @(21, 33),
        @(22, 1),
        @(25, 2915)
    )
    "演出" = ,@(2, 42)
    "全部类型" = @(
        @(2, 1688),
        @(3, 9160),
        @(4, 116),
        @(5, 42),
        @(6, 509),
        @(7, 712),
        @(8, 1380),
        @(9, 207),
        @(11, 101),
        @(12, 5946),
        @(14, 390),
        @(16, 4590),
        @(20, 34),
        @(21, 341),
        @(22, 33),
        @(23, 1),
        @(26, 2915)
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $row = $pair[0]
        $val = $pair[1]
        $ws.Cells.Item($row, 6).Value = $val   # column F = 6
    }
}
